# "issue return tab added"
# Insert two new leading columns ("Weight returned", "Return Date") before the
# existing "Approval Date" column, and append nine new trailing columns after
# the existing "Price Code" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A:J headers two columns to the right by inserting two
# new blank columns at the front of the sheet.
$ws.Range("A:B").Insert()

# Copy the formatting of the (now shifted) existing header cell C1 onto the
# two freshly-inserted header cells so the new headers match the style of
# the rest of the header row.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New leading headers.
$ws.Range("A1").Value = "Weight returned"
$ws.Range("B1").Value = "Return Date"

# New trailing headers appended after the existing "Price Code" column,
# which is now column L following the two-column insert.
$ws.Range("L1").Copy()
$ws.Range("M1:U1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("M1").Value = "Lot Details (If Any)"
$ws.Range("N1").Value = "Original weight"
$ws.Range("O1").Value = "Weight issued"
$ws.Range("P1").Value = "SELECTION YES/NO"
$ws.Range("Q1").Value = "SELECTION CRITERIA"
$ws.Range("R1").Value = "ASKING RATE"
$ws.Range("S1").Value = "BROKER NAME"
$ws.Range("T1").Value = "PARTY NAME"
$ws.Range("U1").Value = "Remarks"
